# This edit performs a cyclic rotation of several body-text
# paragraphs/runs: each location's text is replaced by the text that used
# to live a bit further down the document (the value that falls off the
# end wraps back around to the "Docente" slot). Because several of the
# "new" strings equal some other location's "old" string, a naive global
# Find/Replace pass could cause a later replacement to match text an
# earlier replacement had just inserted. To avoid that, every replacement
# below is scoped to the specific paragraph that should receive it (each
# target string is unique within its own paragraph/the document), so each
# location is touched exactly once using its own untouched original text.

$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $newText) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    # Trim off the trailing paragraph mark so only the run text is
    # replaced (keeps paragraph formatting/pilcrow intact).
    $rng.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1
    $rng.Text = $newText
}

# Objetivos paragraph (previously "Fornecer oportunidade...")
Set-ParagraphText 6 "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."

# Docente(s) Responsável(eis) paragraph (previously "198273 - Domingos Savio Giordani")
Set-ParagraphText 8 "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia Química nos projetos e processos químicos. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."

# Programa resumido paragraph (previously "Plano de Trabalho específico...")
Set-ParagraphText 10 "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia Química.  Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."

# Programa paragraph (previously "Participação do aluno...")
Set-ParagraphText 12 "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

# Avaliação paragraph holds three separate runs (Método / Critério / Norma
# de recuperação) back to back. Replace each one in left-to-right order,
# always re-searching only the portion of the paragraph that comes AFTER
# the end of the previous match, so a freshly-inserted replacement text
# can never be re-matched by a later step, and the still-untouched
# original runs further right are matched unambiguously.
$para14 = $d.Paragraphs.Item(14).Range
$paraEnd = $para14.End
$cursor = $para14.Start

$r = $d.Range($cursor, $paraEnd)
$r.Find.Execute("Supervisão das atividades desenvolvidas pelo aluno durante o estágio.", $true) | Out-Null
$oldLen = $r.End - $r.Start
$r.Text = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."
$paraEnd = $paraEnd + ($r.End - $r.Start) - $oldLen
$cursor = $r.End

$r = $d.Range($cursor, $paraEnd)
$r.Find.Execute("MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio.", $true) | Out-Null
$oldLen = $r.End - $r.Start
$r.Text = "Não será oferecida recuperação."
$paraEnd = $paraEnd + ($r.End - $r.Start) - $oldLen
$cursor = $r.End

$r = $d.Range($cursor, $paraEnd)
$r.Find.Execute("Não será oferecida recuperação.", $true) | Out-Null
$r.Text = "A ser definida com o orientador em função das atividades desenvolvidas no estágio."

# Bibliografia paragraph (previously "A ser definida com o orientador...")
Set-ParagraphText 16 "198273 - Domingos Savio Giordani"
